{"js": "const body = context.document.body;\nconst pairs = [\n  [\"2024-11-06 Wednesday\", \"2024-11-07 Thursday\"],\n  [\"8+52=\", \"92-82=\"],\n  [\"83-20=\", \"86-76=\"],\n  [\"65-21=\", \"0+16=\"],\n  [\"1+85=\", \"69-35=\"],\n  [\"98-38=\", \"83-60=\"],\n  [\"97-1=\", \"90-47=\"],\n  [\"99-16=\", \"44+22=\"],\n  [\"19+58=\", \"53-9=\"],\n  [\"86-16=\", \"99-12=\"],\n  [\"57-2=\", \"8+33=\"],\n  [\"54+17=\", \"87+6=\"],\n  [\"86-57=\", \"38+43=\"],\n  [\"97-91=\", \"96-14=\"],\n  [\"17+49=\", \"92-33=\"],\n  [\"30-29=\", \"19+32=\"],\n  [\"16+19=\", \"36-22=\"],\n  [\"14+6=\", \"45+50=\"],\n  [\"66+27=\", \"97-35=\"],\n  [\"3+31=\", \"27+45=\"],\n  [\"17+42=\", \"83-44=\"],\n  [\"27-7=\", \"94-90=\"],\n  [\"39+42=\", \"46+50=\"],\n  [\"31+66=\", \"6+17=\"],\n  [\"29-24=\", \"39+48=\"],\n  [\"18-14=\", \"88-9=\"],\n  [\"0+2=\", \"91+1=\"],\n  [\"72-50=\", \"57+26=\"],\n  [\"59-37=\", \"17+55=\"],\n  [\"91-45=\", \"15+49=\"],\n  [\"45+43=\", \"26+43=\"],\n  [\"27+55=\", \"79-59=\"],\n  [\"56-39=\", \"15+59=\"],\n  [\"20+39=\", \"75-15=\"],\n  [\"23+62=\", \"10+70=\"],\n  [\"58+17=\", \"87-78=\"],\n  [\"25+45=\", \"28+67=\"],\n  [\"60+11=\", \"67+9=\"],\n  [\"29+60=\", \"79-47=\"],\n  [\"10+8=\", \"80-50=\"],\n  [\"78+7=\", \"15+36=\"],\n  [\"20+68=\", \"3+63=\"],\n  [\"46+26=\", \"98+1=\"],\n  [\"84-59=\", \"58-12=\"],\n  [\"28+0=\", \"10+85=\"],\n  [\"79-58=\", \"44-9=\"],\n  [\"75-0=\", \"67-37=\"],\n  [\"36+7=\", \"36+35=\"],\n  [\"42+32=\", \"92-22=\"],\n  [\"76-63=\", \"31+52=\"],\n  [\"77-29=\", \"43+21=\"],\n  [\"34-27=\", \"28+67=\"],\n  [\"46+37=\", \"67-17=\"],\n  [\"22+69=\", \"62-56=\"],\n  [\"7+32=\", \"23-7=\"],\n  [\"89-9=\", \"45+36=\"],\n  [\"31+55=\", \"46+4=\"],\n  [\"30+32=\", \"28+24=\"],\n  [\"83-12=\", \"90-20=\"],\n  [\"4+43=\", \"80+12=\"],\n  [\"31-4=\", \"0+9=\"],\n  [\"41+16=\", \"26-14=\"],\n  [\"57-55=\", \"8+49=\"],\n  [\"21+38=\", \"18+58=\"],\n  [\"42+4=\", \"31+50=\"],\n  [\"46-36=\", \"59-41=\"],\n  [\"96-34=\", \"58-57=\"],\n  [\"54-9=\", \"77+12=\"],\n  [\"54+36=\", \"11+39=\"],\n  [\"5+41=\", \"62-45=\"],\n  [\"55+17=\", \"46-8=\"],\n  [\"22+32=\", \"80+11=\"],\n  [\"38-16=\", \"52+24=\"],\n  [\"7+38=\", \"10+39=\"],\n  [\"35+52=\", \"95-72=\"],\n  [\"50-33=\", \"47+28=\"],\n  [\"96-35=\", \"44+47=\"],\n  [\"46+31=\", \"62-18=\"],\n  [\"61-12=\", \"72+2=\"],\n  [\"43-19=\", \"44-27=\"],\n  [\"48-4=\", \"13+72=\"],\n  [\"33-26=\", \"55-47=\"],\n  [\"51-48=\", \"88-61=\"],\n  [\"80-72=\", \"7+70=\"],\n  [\"31+8=\", \"83-35=\"],\n  [\"98-14=\", \"2+19=\"],\n  [\"61-14=\", \"22-8=\"],\n  [\"19+22=\", \"3+66=\"],\n  [\"50+28=\", \"37-37=\"],\n  [\"18+8=\", \"46+50=\"],\n  [\"26+11=\", \"2+79=\"],\n  [\"36+23=\", \"89-4=\"],\n  [\"15+45=\", \"76-56=\"],\n  [\"88-88=\", \"55-19=\"],\n  [\"81-42=\", \"4+69=\"],\n  [\"39+27=\", \"27+55=\"],\n  [\"97-53=\", \"28-14=\"],\n  [\"32-13=\", \"5+68=\"],\n  [\"6+78=\", \"85-3=\"],\n  [\"36-17=\", \"9+21=\"],\n  [\"53-44=\", \"9+68=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"NOT FOUND: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n  @('2024-11-06 Wednesday', '2024-11-07 Thursday'),\n  @('8+52=', '92-82='),\n  @('83-20=', '86-76='),\n  @('65-21=', '0+16='),\n  @('1+85=', '69-35='),\n  @('98-38=', '83-60='),\n  @('97-1=', '90-47='),\n  @('99-16=', '44+22='),\n  @('19+58=', '53-9='),\n  @('86-16=', '99-12='),\n  @('57-2=', '8+33='),\n  @('54+17=', '87+6='),\n  @('86-57=', '38+43='),\n  @('97-91=', '96-14='),\n  @('17+49=', '92-33='),\n  @('30-29=', '19+32='),\n  @('16+19=', '36-22='),\n  @('14+6=', '45+50='),\n  @('66+27=', '97-35='),\n  @('3+31=', '27+45='),\n  @('17+42=', '83-44='),\n  @('27-7=', '94-90='),\n  @('39+42=', '46+50='),\n  @('31+66=', '6+17='),\n  @('29-24=', '39+48='),\n  @('18-14=', '88-9='),\n  @('0+2=', '91+1='),\n  @('72-50=', '57+26='),\n  @('59-37=', '17+55='),\n  @('91-45=', '15+49='),\n  @('45+43=', '26+43='),\n  @('27+55=', '79-59='),\n  @('56-39=', '15+59='),\n  @('20+39=', '75-15='),\n  @('23+62=', '10+70='),\n  @('58+17=', '87-78='),\n  @('25+45=', '28+67='),\n  @('60+11=', '67+9='),\n  @('29+60=', '79-47='),\n  @('10+8=', '80-50='),\n  @('78+7=', '15+36='),\n  @('20+68=', '3+63='),\n  @('46+26=', '98+1='),\n  @('84-59=', '58-12='),\n  @('28+0=', '10+85='),\n  @('79-58=', '44-9='),\n  @('75-0=', '67-37='),\n  @('36+7=', '36+35='),\n  @('42+32=', '92-22='),\n  @('76-63=', '31+52='),\n  @('77-29=', '43+21='),\n  @('34-27=', '28+67='),\n  @('46+37=', '67-17='),\n  @('22+69=', '62-56='),\n  @('7+32=', '23-7='),\n  @('89-9=', '45+36='),\n  @('31+55=', '46+4='),\n  @('30+32=', '28+24='),\n  @('83-12=', '90-20='),\n  @('4+43=', '80+12='),\n  @('31-4=', '0+9='),\n  @('41+16=', '26-14='),\n  @('57-55=', '8+49='),\n  @('21+38=', '18+58='),\n  @('42+4=', '31+50='),\n  @('46-36=', '59-41='),\n  @('96-34=', '58-57='),\n  @('54-9=', '77+12='),\n  @('54+36=', '11+39='),\n  @('5+41=', '62-45='),\n  @('55+17=', '46-8='),\n  @('22+32=', '80+11='),\n  @('38-16=', '52+24='),\n  @('7+38=', '10+39='),\n  @('35+52=', '95-72='),\n  @('50-33=', '47+28='),\n  @('96-35=', '44+47='),\n  @('46+31=', '62-18='),\n  @('61-12=', '72+2='),\n  @('43-19=', '44-27='),\n  @('48-4=', '13+72='),\n  @('33-26=', '55-47='),\n  @('51-48=', '88-61='),\n  @('80-72=', '7+70='),\n  @('31+8=', '83-35='),\n  @('98-14=', '2+19='),\n  @('61-14=', '22-8='),\n  @('19+22=', '3+66='),\n  @('50+28=', '37-37='),\n  @('18+8=', '46+50='),\n  @('26+11=', '2+79='),\n  @('36+23=', '89-4='),\n  @('15+45=', '76-56='),\n  @('88-88=', '55-19='),\n  @('81-42=', '4+69='),\n  @('39+27=', '27+55='),\n  @('97-53=', '28-14='),\n  @('32-13=', '5+68='),\n  @('6+78=', '85-3='),\n  @('36-17=', '9+21='),\n  @('53-44=', '9+68='),\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
